# bug 59814: clear evaluation workbook and evaluation sheet caches.
# Add two rows to the "Formulas" sheet that reference the structured
# reference \_Prime.1[Name] (exercises the evaluator caches described in
# the bug), and update the active selections on both sheets.

$wb = $excel.ActiveWorkbook

$wsFormulas = $wb.Worksheets.Item("Formulas")
$wsTable = $wb.Worksheets.Item("Table")

# Add A2 and A3 with the structured reference formula returning the
# "Name" column values ("one" and "two") from the table.
$wsFormulas.Range("A2").Formula = "=\_Prime.1[Name]"
$wsFormulas.Range("A3").Formula = "=\_Prime.1[Name]"

# Update the active cell selection on the "Table" sheet to A7.
$wsTable.Range("A7").Select() | Out-Null

# Update the active cell selection on the "Formulas" sheet to A2, and
# leave that sheet as the selected/active tab.
$wsFormulas.Range("A2").Select() | Out-Null
